# netCrypto.xlsx — "Add files via upload" edit
#
# Semantic changes applied:
#   1. Cell T2 on the active sheet: 647200 -> 0
#   2. Window/view navigation: scroll the sheet so column H is the
#      left-most visible column (was column L), and move the
#      selection/active cell to T3 (was V11).

$wbx = $excel.ActiveWorkbook
$ws  = $wbx.ActiveSheet

# 1. Update the data cell.
$ws.Range("T2").Value = 0

# 2. Scroll the view so the top-left visible cell becomes H1 (column 8,
#    row 1), matching a user having scrolled left from column L.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1

# Reposition the window itself to match the saved workbook view state.
$win.Left = -120
$win.Top = -120

# 3. Move the selection/active cell to T3.
$ws.Range("T3").Select()
